$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 245; this shifts the existing rows
# 245-368 down to 246-369 (and grows the used range / dimension
# accordingly), matching the author's edit of adding one new weekly
# record ahead of the previously-first record in this block.
$ws.Rows("245:245").Insert()

# Populate the newly inserted row 245 with the new weekly price record.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R share the same constant values as every
# other row in this "Hortaliza, Terminal La Palmera de La Serena -
# Espinaca" sheet.
$ws.Range("A245").Value = 8
$ws.Range("B245").Value = "Terminal La Palmera de La Serena"
$ws.Range("C245").Value = "Coquimbo"
$ws.Range("D245").Value = 44960
$ws.Range("E245").Value = 4
$ws.Range("F245").Value = 100112012
$ws.Range("G245").Value = "Espinaca"
$ws.Range("H245").Value = "Sin especificar"
$ws.Range("I245").Value = "Primera"
$ws.Range("J245").Value = 1360
$ws.Range("K245").Value = 500
$ws.Range("L245").Value = 600
$ws.Range("M245").Value = 550
$ws.Range("N245").Value = "$/atado 300 a 500 gramos"
$ws.Range("O245").Value = "Provincia del Elquí"
$ws.Range("P245").Value = 1100
$ws.Range("Q245").Value = 0.5
$ws.Range("R245").Value = "Hortaliza"
